$wb = $excel.ActiveWorkbook

# --- Sheet "БИВТ-22-17" (grading sheet 1) ---
$ws1 = $wb.Worksheets.Item("БИВТ-22-17")
$ws1.Activate()

# Row 13 previously had a stray "pass" text in the 3rd-lab column; correct it
# to the actual numeric score (the dependent Score / Mid-term formulas will
# recalc automatically).
$ws1.Range("D13").Value = 5

$ws1.Range("B38").Select()

# --- Sheet "БИВТ-22-20" (grading sheet 3) ---
$ws3 = $wb.Worksheets.Item("БИВТ-22-20")
$ws3.Activate()

$ws3.Range("F3").Value = 5
$ws3.Range("F8").Value = 5

# F12 carries a "number stored as text" quote-prefix style (it used to hold
# the stray "pass" label); writing a plain numeric .Value resets that style
# flag, so round-trip the formatting through a scratch cell to keep the
# original style intact while still landing the corrected score.
$scratch3 = $ws3.Range("ZZ500")
$ws3.Range("F12").Copy()
$scratch3.PasteSpecial(-4122)
$ws3.Range("F12").Value = 4
$scratch3.Copy()
$ws3.Range("F12").PasteSpecial(-4122)
$scratch3.Clear()
$excel.CutCopyMode = $false

$ws3.Range("E20").Select()

# --- Sheet "БИВТ-22-18" (grading sheet 2) ---
$ws2 = $wb.Worksheets.Item("БИВТ-22-18")
$ws2.Activate()

$ws2.Range("F20").Value = 5
$ws2.Range("F21").Value = 5
$ws2.Range("E26").Value = 5

# J26 previously held a one-off formula that (incorrectly) referenced H32
# instead of H26, breaking the fill pattern used by J4:J25. Bring it back in
# line with the rest of the column, matching both formula and formatting of
# the cell above it.
$ws2.Range("J25").Copy()
$ws2.Range("J26").PasteSpecial(-4122)
$ws2.Range("J26").Formula = "=IF(H26>=25,5,IF(H26>=22,4,IF(H26>=19,3,2)))"

$ws2.Range("F21").Select()

$excel.CutCopyMode = $false
